# Generate Report for Handback
# - Both locale sheets (zh-cn, de-de) moved from "Ready for handoff" to
#   "Handed back: in sync with en-US" once the handback round-tripped cleanly.
# - The Overview roll-up sheet mirrors that same status text per locale column.
# - The "Latest Handback DateTime" for each locale is refreshed to the
#   timestamp of this run.
# - The stale "handback file is not the latest" warning is cleared now that
#   the handback is in sync.
# - Column widths are re-fitted to the new cell content (status text grew,
#   the error-detail text shrank to nothing).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Status column on the per-locale sheets
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Overview sheet mirrors the same status for each locale column
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Latest Handback DateTime refreshed
$wsZhCn.Range("K2").Value = "2016-08-17 06:43:29"
$wsDeDe.Range("K2").Value = "2016-08-17 06:43:37"

# Error Detail no longer applies now the handback is in sync - clear it
$wsZhCn.Range("P2").ClearContents()
$wsDeDe.Range("P2").ClearContents()

# Re-fit the columns whose content just changed
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()
$wsZhCn.Columns.Item(3).AutoFit()
$wsZhCn.Columns.Item(16).AutoFit()
$wsDeDe.Columns.Item(3).AutoFit()
$wsDeDe.Columns.Item(16).AutoFit()
